# "clean up data_summary_table for compare/contrast, brute force add radius differences"
#
# The summary table on Sheet1 (site gps.xlsx) is refreshed with recomputed
# per-site averages/extremes (columns D:J, rows 2:4), and the sheet's active
# selection is moved to J2 after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Alegria (row 2) ---
$ws.Range("D2").Value = 7.98
$ws.Range("E2").Value = 17.27
$ws.Range("G2").Value = 8.23
$ws.Range("H2").Value = 7.63
$ws.Range("J2").Value = 14.06

# --- Bodega Bay (row 3) ---
$ws.Range("D3").Value = 7.89
$ws.Range("E3").Value = 13.92
$ws.Range("G3").Value = 8.27
$ws.Range("H3").Value = 7.5
$ws.Range("I3").Value = 17.18
$ws.Range("J3").Value = 10.01

# --- Lompoc Landing (row 4) ---
$ws.Range("D4").Value = 7.93
$ws.Range("E4").Value = 14.36
$ws.Range("G4").Value = 8.25
$ws.Range("H4").Value = 7.28
$ws.Range("J4").Value = 11.9

# Scroll the viewport so column D is left-most visible (best effort — some
# hosts only persist this alongside frozen/split panes) and land the active
# selection on J2, matching the author's final cursor position.
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("J2").Select()
